$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data columns were re-ordered: the "4:1 Conditioned" /
# "1:4 Unconditioned" columns (previously D/E) move to B/D, and the
# "1:4 Conditioned" / "4:1 Unconditioned" columns (previously B/C) move to C/E.
# Concretely: new B = old D, new C = old B, new D = old E, new E = old C.
for ($r = 1; $r -le 9; $r++) {
    $oldB = $ws.Cells.Item($r, 2).Value2
    $oldC = $ws.Cells.Item($r, 3).Value2
    $oldD = $ws.Cells.Item($r, 4).Value2
    $oldE = $ws.Cells.Item($r, 5).Value2

    $ws.Cells.Item($r, 2).Value2 = $oldD
    $ws.Cells.Item($r, 3).Value2 = $oldB
    $ws.Cells.Item($r, 4).Value2 = $oldE
    $ws.Cells.Item($r, 5).Value2 = $oldC
}

# Columns D and E now share the same width (previously D=15, E=15.83203125).
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# Update the view: zoom in to 125% and move the selection to F7.
$excel.ActiveWindow.Zoom = 125
$ws.Range("F7").Select() | Out-Null
